$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 1.02
$ws.Range("K2").Value = 19
$ws.Range("N2").Value = 1.44
$ws.Range("O2").Value = 2.75
$ws.Range("W2").Value = 23
$ws.Range("AI2").Value = 23
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 7.5
$ws.Range("J3").Value = 1.14
$ws.Range("K3").Value = 5.5
$ws.Range("R3").Value = 2.75
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 4.33
$ws.Range("AB3").Value = 29
$ws.Range("AH3").Value = 81
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 3.2
$ws.Range("L4").Value = 1.57
$ws.Range("M4").Value = 2.25
$ws.Range("N4").Value = 2.88
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.62
$ws.Range("Q4").Value = 2.2
$ws.Range("T4").Value = 6
$ws.Range("U4").Value = 10
$ws.Range("V4").Value = 11
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 26
$ws.Range("AE4").Value = 7
$ws.Range("AF4").Value = 13
$ws.Range("AH4").Value = 34
$ws.Range("AJ4").Value = 41
$ws.Range("G6").Value = 2.75
$ws.Range("I6").Value = 2.75
$ws.Range("X6").Value = 34
$ws.Range("AF6").Value = 11
$ws.Range("G9").Value = 2.15
$ws.Range("I9").Value = 3.3
$ws.Range("K9").Value = 12
$ws.Range("T9").Value = 9.5
$ws.Range("U9").Value = 12
$ws.Range("W9").Value = 21
$ws.Range("X9").Value = 17
$ws.Range("AE9").Value = 12
$ws.Range("AH9").Value = 34
$ws.Range("AI9").Value = 23
$ws.Range("N10").Value = 1.95
$ws.Range("O10").Value = 1.85
$ws.Range("K11").Value = 7.5
$ws.Range("L11").Value = 1.36
$ws.Range("M11").Value = 3
$ws.Range("G13").Value = 1.5
$ws.Range("H13").Value = 4.2
$ws.Range("I13").Value = 5.75
$ws.Range("J13").Value = 1.06
$ws.Range("K13").Value = 10
$ws.Range("L13").Value = 1.33
$ws.Range("M13").Value = 3.25
$ws.Range("N13").Value = 2.08
$ws.Range("O13").Value = 1.73
$ws.Range("R13").Value = 2.2
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 6
$ws.Range("U13").Value = 6.5
$ws.Range("V13").Value = 8.5
$ws.Range("W13").Value = 10
$ws.Range("AA13").Value = 8.5
$ws.Range("AF13").Value = 29
$ws.Range("AG13").Value = 19
$ws.Range("AH13").Value = 67
$ws.Range("N15").Value = 1.73
$ws.Range("O15").Value = 2.08
$ws.Range("G17").Value = 2.82
$ws.Range("H17").Value = 3.85
$ws.Range("I17").Value = 2.15
$ws.Range("L17").Value = 1.13
$ws.Range("M17").Value = 5.1
$ws.Range("Q17").Value = 3.65
$ws.Range("S17").Value = 2.75
$ws.Range("T17").Value = 15.5
$ws.Range("U17").Value = 19.5
$ws.Range("V17").Value = 10.75
$ws.Range("W17").Value = 37
$ws.Range("X17").Value = 20
$ws.Range("Y17").Value = 21
$ws.Range("AE17").Value = 14
$ws.Range("AF17").Value = 15
$ws.Range("AG17").Value = 9.25
$ws.Range("AH17").Value = 23
$ws.Range("AJ17").Value = 17.5
$ws.Range("G18").Value = 2.05
$ws.Range("I18").Value = 3
$ws.Range("L18").Value = 1.11
$ws.Range("M18").Value = 5.6
$ws.Range("N18").Value = 1.36
$ws.Range("O18").Value = 2.92
$ws.Range("P18").Value = 1.21
$ws.Range("Q18").Value = 3.9
$ws.Range("R18").Value = 1.34
$ws.Range("S18").Value = 3
$ws.Range("T18").Value = 15
$ws.Range("U18").Value = 15.5
$ws.Range("V18").Value = 9.25
$ws.Range("W18").Value = 23
$ws.Range("Z18").Value = 10.25
$ws.Range("AA18").Value = 9
$ws.Range("AB18").Value = 10.5
$ws.Range("AC18").Value = 26
$ws.Range("AD18").Value = 110
$ws.Range("AE18").Value = 19
$ws.Range("AF18").Value = 24
$ws.Range("AG18").Value = 11.5
$ws.Range("AI18").Value = 21
$ws.Range("AJ18").Value = 19.5
$ws.Range("G19").Value = 3.4
$ws.Range("I19").Value = 1.93
$ws.Range("L19").Value = 1.15
$ws.Range("N19").Value = 1.47
$ws.Range("O19").Value = 2.5
$ws.Range("P19").Value = 1.25
$ws.Range("Q19").Value = 3.55
$ws.Range("T19").Value = 16.5
$ws.Range("U19").Value = 24
$ws.Range("X19").Value = 25
$ws.Range("AA19").Value = 8
$ws.Range("AH19").Value = 18.5
$ws.Range("G21").Value = 1.8
$ws.Range("I21").Value = 4.33
$ws.Range("J21").Value = 1.07
$ws.Range("K21").Value = 8.5
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 1.75
$ws.Range("W21").Value = 15
$ws.Range("X21").Value = 17
$ws.Range("AE21").Value = 10
$ws.Range("AF21").Value = 21
$ws.Range("H27").Value = 6
$ws.Range("I27").Value = 17
$ws.Range("L27").Value = 1.18
$ws.Range("M27").Value = 4.5
$ws.Range("N27").Value = 1.6
$ws.Range("O27").Value = 2.3
$ws.Range("R27").Value = 2.38
$ws.Range("S27").Value = 1.53
$ws.Range("AA27").Value = 12
$ws.Range("AB27").Value = 29
$ws.Range("AF27").Value = 67
$ws.Range("N47").Value = 2.88
$ws.Range("O47").Value = 1.4
$ws.Range("AI17").Value = 14.5
$ws.Range("X18").Value = 14
